$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 74-77 (dates already present, only MarketObjects column changes)
$ws.Range("B74").Value = "['BTCUSD.SPOT']"
$ws.Range("B75").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"
$ws.Range("B76").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"
$ws.Range("B77").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"

# Row 78 already has its date; set the MarketObjects column
$ws.Range("B78").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']"

# New rows 79-87: Date column + MarketObjects column.
# Dates must be written as text, not auto-converted to date serials, so the
# target cell is temporarily formatted as Text, then restored to the default
# "Normal" style so no explicit style index is left behind on the cell.
$newRows = @(
    @{ Row = 79; Date = "2025-08-29"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']" },
    @{ Row = 80; Date = "2025-08-30"; Objs = "['BTCUSD.SPOT']" },
    @{ Row = 81; Date = "2025-08-31"; Objs = "['BTCUSD.SPOT']" },
    @{ Row = 82; Date = "2025-09-01"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']" },
    @{ Row = 83; Date = "2025-09-02"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']" },
    @{ Row = 84; Date = "2025-09-03"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']" },
    @{ Row = 85; Date = "2025-09-04"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE.REGULARIZED']" },
    @{ Row = 86; Date = "2025-09-05"; Objs = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']" },
    @{ Row = 87; Date = "2025-09-06"; Objs = "[]" }
)

foreach ($r in $newRows) {
    $dateCell = $ws.Cells.Item($r.Row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.Date
    $dateCell.Style = "Normal"

    $objCell = $ws.Cells.Item($r.Row, 2)
    $objCell.Value = $r.Objs
}
